$d = $word.ActiveDocument

# --- Part 1: insert a new "Meta description" paragraph right after the title (paragraph 1) ---
# Do this first, while the trailing duplicate-title paragraph (which has the exact
# leading-empty-run + bold-run shape we want to reuse) still exists at the end.
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

# Borrow the run/formatting layout (leading empty run + bold run) from the
# duplicate-title paragraph that still lives at the end of the document, so the new
# paragraph's XML shape (<w:r/> + bold run) matches the rest of the doc's idiom.
$srcCount = $d.Paragraphs.Count
$boldShapeSrc = $d.Paragraphs($srcCount - 1)
$metaPara.Range.FormattedText = $boldShapeSrc.Range.FormattedText

$metaPara = $d.Paragraphs(2)
$boldTextStart = $metaPara.Range.Start
$boldTextEnd = $metaPara.Range.End - 1
$boldRange = $d.Range($boldTextStart, $boldTextEnd)
$boldRange.Text = "Meta description"

$metaPara = $d.Paragraphs(2)
$tailPos = $metaPara.Range.End - 1
$tail = $d.Range($tailPos, $tailPos)
$tail.InsertAfter(": Enter the eerie world of Black Hawk Deluxe offering customizable user interface and Volatility Levels for winning opportunities. Play now for free.")

$metaPara = $d.Paragraphs(2)
$plainStart = $tailPos
$plainEnd = $metaPara.Range.End - 1
$plainRange = $d.Range($plainStart, $plainEnd)
$plainRange.Bold = 0

# --- Part 2: remove the trailing bold "Play Black Hawk Deluxe Free | ..." paragraph ---
# It's the penultimate paragraph (the last one is the italic meta-description paragraph).
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($count - 1)
$dupTitlePara.Range.Delete()

# --- Part 3: rewrite the final (italic) paragraph's text with the new image prompt ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastRange = $lastPara.Range
$lastTextRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$lastTextRange.Text = "Create a cartoon-style feature image for ""Black Hawk Deluxe"" featuring a happy Maya warrior with glasses. The Maya warrior should be depicted as brave, ready to face the legion of undead soldiers to claim the treasure of Black Hawk Castle. The background of the image should be the menacing castle looming over the reels, with the undead soldiers shuffling towards it. The overall style should be dark and eerie but with a touch of humor, capturing the adventurous spirit of the game. The image should be eye-catching and appealing to players who enjoy horror and fantasy-themed slots."
